# Bazdid > report 6-2: add a new "Parts filter" report-date row (row 22)
# to the dates table, mirroring the format of the row above it (row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 22 -----------------------------------------------------
# Copy the formatting (fill/border/style) of the row above (row 20)
# onto the new row so A22:C22 match the look of the other data rows.
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null

# C22 ("1400/04/26") looks like a date and would otherwise be silently
# auto-converted to a date serial on assignment; entering it as a
# formula that evaluates to the literal text, then collapsing the
# formula down to its value, keeps it a plain text cell.
$ws.Range("C22").Formula = "=""1400/04/26"""
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null

$ws.Range("B22").Value = "2021 Jully 17"
$ws.Range("A22").Value = "Report 6-2 (Bazdid) *"

# --- View state -------------------------------------------------------
# Scroll the sheet so row 8 is at the top and select the new row's
# first cell, matching the author's on-screen state after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("A22").Select() | Out-Null
